$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Step1 - Input Data")
$ws2 = $wb.Worksheets.Item("Step2 - Projection")

# Company Ticker Symbol: mmm -> MMM
$ws1.Range("B3").Value = "MMM"

# Personal Required Rate of Return: 5.87% -> 5.88%
# (stored as literal text, not a numeric percentage, so force the
# number format back to its original percent format afterwards to
# avoid leaving a "quote prefix" / text-format style behind)
$ws1.Range("B4").NumberFormat = "@"
$ws1.Range("B4").Value = "5.88%"
$ws1.Range("B4").NumberFormat = "0%"

# Growth Rate: 0.45% -> 0.5%
$ws2.Range("C5").NumberFormat = "@"
$ws2.Range("C5").Value = "0.5%"
$ws2.Range("C5").NumberFormat = "0%"

$wb.Application.CalculateFullRebuild()
